# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit refresh to Sheets/Asura_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 474.5
$ws.Range("I2").Value = 450
$ws.Range("K2").Value = 450
$ws.Range("M2").Value = -337

# Row 6
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 240
$ws.Range("K6").Value = 720
$ws.Range("M6").Value = -608

# Row 9
$ws.Range("H9").Value = 109.166664
$ws.Range("I9").Value = 109.875
$ws.Range("J9").Value = 107.75
$ws.Range("K9").Value = 109.875
$ws.Range("L9").Value = 107.75
$ws.Range("M9").Value = 59.125
$ws.Range("N9").Value = -445.75

# Row 12
$ws.Range("H12").Value = 2002.2
$ws.Range("I12").Value = 2003.6666
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 2003.6666
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = -1833.6666
$ws.Range("N12").Value = -2340

# Row 21
$ws.Range("H21").Value = 26000
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20936

# Row 23
$ws.Range("H23").Value = 26000
$ws.Range("J23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("N23").Value = -20468

# Row 28
$ws.Range("H28").Value = 564.63336
$ws.Range("I28").Value = 563.6786
$ws.Range("K28").Value = 563.6786
$ws.Range("M28").Value = -78.67859999999996

# Row 29
$ws.Range("H29").Value = 24
$ws.Range("I29").Value = 24
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 72
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 209
$ws.Range("N29").ClearContents()

# Row 38
$ws.Range("H38").Value = 2086
$ws.Range("I38").Value = 108.125
$ws.Range("J38").Value = 9997.5
$ws.Range("K38").Value = 324.375
$ws.Range("L38").Value = 29992.5
$ws.Range("M38").Value = 47.625
$ws.Range("N38").Value = -30736.5

# Row 51
$ws.Range("H51").Value = 2499.75
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 2999.5
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 2999.5
$ws.Range("M51").Value = -1516
$ws.Range("N51").Value = -3967.5

# Row 58
$ws.Range("H58").Value = 1753.6666
$ws.Range("I58").Value = 769
$ws.Range("J58").Value = 2541.4
$ws.Range("K58").Value = 2307
$ws.Range("L58").Value = 7624.200000000001
$ws.Range("M58").Value = -2157
$ws.Range("N58").Value = -7924.200000000001

# Row 62
$ws.Range("H62").Value = 1698.75
$ws.Range("I62").Value = 1448.3334
$ws.Range("K62").Value = 1448.3334
$ws.Range("M62").Value = -824.3334

# Row 65
$ws.Range("H65").Value = 1698.75
$ws.Range("I65").Value = 1448.3334
$ws.Range("K65").Value = 7241.666999999999
$ws.Range("M65").Value = -4121.666999999999

# Row 92
$ws.Range("H92").Value = 2276

# Row 107
$ws.Range("H107").Value = 843.3333
$ws.Range("I107").Value = 843.3333
$ws.Range("K107").Value = 843.3333
$ws.Range("M107").Value = 1076.6667

# Row 112
$ws.Range("H112").Value = 3432.2666
$ws.Range("J112").Value = 3606
$ws.Range("L112").Value = 10818
$ws.Range("N112").Value = -13034

# Row 129
$ws.Range("H129").Value = 953.72
$ws.Range("J129").Value = 972.0103
$ws.Range("L129").Value = 2916.0309
$ws.Range("N129").Value = -12916.0309

# Row 131
$ws.Range("H131").Value = 4249.7744
$ws.Range("J131").Value = 8667.154
$ws.Range("L131").Value = 26001.462
$ws.Range("N131").Value = -36081.462

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10071.912
$ws.Range("I32").Value = 10762.151
$ws.Range("K32").Value = 10762.151
$ws.Range("M32").Value = -10475.151

# Row 37
$ws.Range("H37").Value = 12593.333
$ws.Range("J37").Value = 12593.333
$ws.Range("L37").Value = 12593.333
$ws.Range("N37").Value = -13139.333

# Row 44
$ws.Range("H44").Value = 34450
$ws.Range("J44").Value = 34450
$ws.Range("L44").Value = 34450
$ws.Range("N44").Value = -35426

# Row 61
$ws.Range("H61").Value = 2237.6875
$ws.Range("I61").Value = 2037.75
$ws.Range("J61").Value = 2837.5
$ws.Range("K61").Value = 2037.75
$ws.Range("L61").Value = 2837.5
$ws.Range("M61").Value = -1825.75
$ws.Range("N61").Value = -3261.5

# Row 113
$ws.Range("H113").Value = 40000
$ws.Range("J113").Value = 40000
$ws.Range("L113").Value = 40000
$ws.Range("N113").Value = -48678

# Row 132
$ws.Range("H132").Value = 28001.25
$ws.Range("I132").Value = 56506
$ws.Range("K132").Value = 169518
$ws.Range("M132").Value = -166988

# Row 136
$ws.Range("H136").Value = 2237.6875
$ws.Range("I136").Value = 2037.75
$ws.Range("J136").Value = 2837.5
$ws.Range("K136").Value = 6113.25
$ws.Range("L136").Value = 8512.5
$ws.Range("M136").Value = -3563.25
$ws.Range("N136").Value = -13612.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 122
$ws.Range("H122").Value = 42780
$ws.Range("J122").Value = 42780
$ws.Range("L122").Value = 42780
$ws.Range("N122").Value = -52580

# Row 132
$ws.Range("H132").Value = 94289.09
$ws.Range("J132").Value = 94289.09
$ws.Range("L132").Value = 94289.09
$ws.Range("N132").Value = -104409.09

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1831.4464
$ws.Range("I31").Value = 2325.3704
$ws.Range("J31").Value = 1371.5862
$ws.Range("K31").Value = 2325.3704
$ws.Range("L31").Value = 1371.5862
$ws.Range("M31").Value = -2030.3704
$ws.Range("N31").Value = -1961.5862

# Row 34
$ws.Range("H34").Value = 1831.4464
$ws.Range("I34").Value = 2325.3704
$ws.Range("J34").Value = 1371.5862
$ws.Range("K34").Value = 2325.3704
$ws.Range("L34").Value = 1371.5862
$ws.Range("M34").Value = -2123.3704
$ws.Range("N34").Value = -1775.5862

# Row 39
$ws.Range("H39").Value = 7514.5557
$ws.Range("I39").Value = 8266.375
$ws.Range("J39").Value = 1500
$ws.Range("K39").Value = 8266.375
$ws.Range("L39").Value = 1500
$ws.Range("M39").Value = -7875.375
$ws.Range("N39").Value = -2282

# Row 49
$ws.Range("H49").Value = 7514.5557
$ws.Range("I49").Value = 8266.375
$ws.Range("J49").Value = 1500
$ws.Range("K49").Value = 8266.375
$ws.Range("L49").Value = 1500
$ws.Range("M49").Value = -8084.375
$ws.Range("N49").Value = -1864

# Row 99
$ws.Range("H99").Value = 1716.909
$ws.Range("I99").Value = 1717.3334
$ws.Range("J99").Value = 1715
$ws.Range("K99").Value = 1717.3334
$ws.Range("L99").Value = 1715
$ws.Range("M99").Value = -219.3334
$ws.Range("N99").Value = -4711

# Row 107
$ws.Range("H107").Value = 661
$ws.Range("I107").Value = 680
$ws.Range("J107").Value = 639.2857
$ws.Range("K107").Value = 680
$ws.Range("L107").Value = 639.2857
$ws.Range("M107").Value = 1240
$ws.Range("N107").Value = -4479.2857

# Row 126
$ws.Range("H126").Value = 1716.909
$ws.Range("I126").Value = 1717.3334
$ws.Range("J126").Value = 1715
$ws.Range("K126").Value = 5152.0002
$ws.Range("L126").Value = 5145
$ws.Range("M126").Value = -2682.0002
$ws.Range("N126").Value = -10085

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 700
$ws.Range("I92").Value = 300
$ws.Range("J92").Value = 900
$ws.Range("K92").Value = 900
$ws.Range("L92").Value = 2700
$ws.Range("M92").Value = 348
$ws.Range("N92").Value = -5196

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1144.32
$ws.Range("I97").Value = 1194.2142
$ws.Range("J97").Value = 1080.8182
$ws.Range("K97").Value = 1194.2142
$ws.Range("L97").Value = 1080.8182
$ws.Range("M97").Value = -698.2141999999999
$ws.Range("N97").Value = -2072.8182

# Row 133
$ws.Range("H133").Value = 71097.5
$ws.Range("J133").Value = 71097.5
$ws.Range("L133").Value = 71097.5
$ws.Range("N133").Value = -81217.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1348.3334
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1418
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1418
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -2008

# Row 27
$ws.Range("H27").Value = 1348.3334
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1418
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1418
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1632

# Row 38
$ws.Range("H38").Value = 50000
$ws.Range("I38").Value = 50000
$ws.Range("K38").Value = 50000
$ws.Range("M38").Value = -49590

# Row 136
$ws.Range("H136").Value = 1389.3334
$ws.Range("I136").Value = 996.275
$ws.Range("J136").Value = 2818.6365
$ws.Range("K136").Value = 2988.825
$ws.Range("L136").Value = 8455.9095
$ws.Range("M136").Value = -438.8249999999998
$ws.Range("N136").Value = -13555.9095

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 95
$ws.Range("H95").Value = 114482.29
$ws.Range("J95").Value = 114482.29
$ws.Range("L95").Value = 114482.29
$ws.Range("N95").Value = -119974.29

# Row 107
$ws.Range("H107").Value = 599.125
$ws.Range("I107").Value = 364.66666
$ws.Range("J107").Value = 739.8
$ws.Range("K107").Value = 1093.99998
$ws.Range("L107").Value = 2219.4
$ws.Range("M107").Value = 826.0000199999999
$ws.Range("N107").Value = -6059.4

# Row 119
$ws.Range("H119").Value = 275349
$ws.Range("J119").Value = 275349
$ws.Range("L119").Value = 275349
$ws.Range("N119").Value = -285025

# Row 132
$ws.Range("H132").Value = 2155.3157
$ws.Range("I132").Value = 1997.76
$ws.Range("K132").Value = 5993.28
$ws.Range("M132").Value = -3463.28

# Row 136
$ws.Range("H136").Value = 2364.7368
$ws.Range("I136").Value = 2397.8096
$ws.Range("J136").Value = 2323.8823
$ws.Range("K136").Value = 7193.4288
$ws.Range("L136").Value = 6971.646900000001
$ws.Range("M136").Value = -4643.4288
$ws.Range("N136").Value = -12071.6469
